# Rename the inline picture shapes.
#
#   footer1 (Pearson logo) : image2.png -> image1.png
#   footer2 (Pearson logo) : image2.png -> image1.png
#   header  (BTEC logo)    : image1.jpg -> image2.jpg
#
$d = $word.ActiveDocument

for ($i = 1; $i -le 2; $i++) {
    $footer = $d.Sections(1).Footers.Item($i)
    if ($footer.Exists -and $footer.Range.InlineShapes.Count -gt 0) {
        $shape = $footer.Range.InlineShapes.Item(1)
        if ($shape.AlternativeText -like "*PearsonLogo.png") {
            $shape.Name = "image1.png"
        }
    }
}

for ($i = 1; $i -le 3; $i++) {
    $header = $d.Sections(1).Headers.Item($i)
    if ($header.Exists -and $header.Range.InlineShapes.Count -gt 0) {
        $shape = $header.Range.InlineShapes.Item(1)
        if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
            $shape.Name = "image2.jpg"
        }
    }
}
